$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values (row 3 & 4 gain an E column response-code entry) ---
$ws.Range("E3").Value = "REFRESH_LANGUAGE_PUSHER"
$ws.Range("E4").Value = "NOTIFY_USER"

# --- New rows 9-11 (row 9 previously only had D9 placeholder) ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Add language"
$ws.Range("C9").Value = "RSNT_GUEST_DEV_<<orgId>>"
$ws.Range("E9").Value = "ADD_LANGUAGE_PUSHER"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "delete language"
$ws.Range("C10").Value = "RSNT_GUEST_DEV_<<orgId>>"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Update Organization setting"
$ws.Range("C11").Value = "RSNT_GUEST_DEV_<<orgId>>"
$ws.Range("E11").Value = "ORG_SETTING_PUSHER"

# --- Formatting: reuse existing bordered-cell looks via copy/paste-format ---
$ws.Range("A3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("A9:A11").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("C9:C11").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("B10:B11").PasteSpecial(-4122)

$ws.Range("E6").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("E7:E8").Borders.Item(7).LineStyle = 1
$ws.Range("E7:E8").Borders.Item(10).LineStyle = 1
$ws.Range("E7:E8").Borders.Item(8).LineStyle = 1

$ws.Range("E9").Borders.Item(7).LineStyle = 1
$ws.Range("E9").Borders.Item(10).LineStyle = 1
$ws.Range("E9").Borders.Item(8).LineStyle = 1

$ws.Range("E11").Borders.Item(7).LineStyle = 1
$ws.Range("E11").Borders.Item(10).LineStyle = 1
$ws.Range("E11").Borders.Item(9).LineStyle = 1

$ws.Application.CutCopyMode = 0

# --- View state ---
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("D30").Select()

Write-Host "done"
